# "Generate Report for Handoff" - refresh the localization-status report
# with a new handoff run: new source GUID, new xliff hash, refreshed
# timestamps, and a blanked-out "Latest Target File" column (no target
# file produced yet for this run).

$wb = $excel.ActiveWorkbook

$oldGuid = "a5c131ed-efda-418f-80f6-33a7da52d47e"
$newGuid = "9d746d88-19a6-4765-b540-e0825f50c8be"
$oldHash = "fe6262f96bb8ccd23f1885eec2bb6c73eb816dfe"
$newHash = "3e4fcabbc87c2d113f5f9681a15d3a34545d538e"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"
$oldMdPath = "e2e\$oldGuid.md"
$newMdPath = "e2e\$newGuid.md"

$newHoDate = "2016-09-01 09:17:13"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"
$newZhHandoffDate = "2016-09-01 09:17:00"
$emptyHandbackDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b85b249c654ab7b21d56fdcb386aae33836c538e/e2e/$newGuid.md", "", "", $newMdPath)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Font.Underline = $false
$wsZh.Range("I2").Font.ColorIndex = -4105
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $emptyHandbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b85b249c654ab7b21d56fdcb386aae33836c538e/e2e/$newGuid.md", "", "", $newMd)

$wsZh.Columns.Item(9).ColumnWidth = 17.8333
$wsZh.Columns.Item(10).ColumnWidth = 20.85

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Font.Underline = $false
$wsDe.Range("I2").Font.ColorIndex = -4105
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $emptyHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b85b249c654ab7b21d56fdcb386aae33836c538e/e2e/$newGuid.md", "", "", $newMd)

$wsDe.Columns.Item(9).ColumnWidth = 17.8333
$wsDe.Columns.Item(10).ColumnWidth = 20.85
